# Generate Report for Handoff
# Updates the localization-status workbook:
#  - refreshes the "md" handoff entry (new guid / new handoff timestamps / new xlf hash)
#  - appends two new rows (the two .png dependency files) on every sheet

$wb = $excel.ActiveWorkbook

$oldGuid = "1ef04d2e-cc4f-4d2d-b3ee-49c311c0b4f7"
$newGuid = "414cc058-a02a-4154-bf74-11e5db4ce015"

$newMdName   = "$newGuid.md"
$newHash     = "673c686a6d3fe4e88443de80c7a94caeb4280e85"
$png1        = "a0bd12a1-03cb-4762-a795-3ad3fb4a1e91.png"
$png2        = "f92b2024-fb18-4693-a9a1-a7964050eb55.png"
$png1target  = "3aa087b4579bb1cdab7d10e05f9c9ff687a8b479.png"
$png2target  = "c7090b485198f949b4d93a599b145b6e300d44e4.png"

$handoffDate      = "2016-03-21 21:00:24"
$handoffDatetime  = "2016-03-21 21:00:13"
$epoch            = "0001-01-01 00:00:00"
$ready            = "Ready for handoff"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $handoffDate

$wsOverview.Rows.Item(3).Insert()
$wsOverview.Range("A3").Value = $png1
$wsOverview.Range("B3").Value = $ready
$wsOverview.Range("C3").Value = $ready
$wsOverview.Range("D3").Value = $handoffDate

$wsOverview.Rows.Item(4).Insert()
$wsOverview.Range("A4").Value = $png2
$wsOverview.Range("B4").Value = $ready
$wsOverview.Range("C4").Value = $ready
$wsOverview.Range("D4").Value = $handoffDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/be33dc629c1990b388306ac50dc5ef87532c4651/e2e/$png1", "", "", $png1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/be33dc629c1990b388306ac50dc5ef87532c4651/e2e/$png2", "", "", $png2)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhXlf = "$newGuid.$newHash.zh-cn.xlf"

$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $zhXlf
$wsZh.Range("E2").Value = $handoffDatetime

$wsZh.Rows.Item(3).Insert()
$wsZh.Range("A3").Value = $png1
$wsZh.Range("B3").Value = ".png"
$wsZh.Range("C3").Value = $ready
$wsZh.Range("D3").Value = $png1target
$wsZh.Range("E3").Value = $handoffDatetime
$wsZh.Range("H3").Value = $epoch
$wsZh.Range("J3").Value = "IsDependency"
$wsZh.Range("K3").Value = "e2e\$newMdName"

$wsZh.Rows.Item(4).Insert()
$wsZh.Range("A4").Value = $png2
$wsZh.Range("B4").Value = ".png"
$wsZh.Range("C4").Value = $ready
$wsZh.Range("D4").Value = $png2target
$wsZh.Range("E4").Value = $handoffDatetime
$wsZh.Range("H4").Value = $epoch
$wsZh.Range("J4").Value = "IsDependency"
$wsZh.Range("K4").Value = "e2e\$newMdName"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/be33dc629c1990b388306ac50dc5ef87532c4651/e2e/$png1", "", "", $png1)
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46cbf783add8b88114af48b7e12de4ed763b940a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png1target", "", "", $png1target)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/be33dc629c1990b388306ac50dc5ef87532c4651/e2e/$png2", "", "", $png2)
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46cbf783add8b88114af48b7e12de4ed763b940a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png2target", "", "", $png2target)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$deXlf = "$newGuid.$newHash.de-de.xlf"

$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $deXlf
$wsDe.Range("E2").Value = $handoffDate

$wsDe.Rows.Item(3).Insert()
$wsDe.Range("A3").Value = $png1
$wsDe.Range("B3").Value = ".png"
$wsDe.Range("C3").Value = $ready
$wsDe.Range("D3").Value = $png1target
$wsDe.Range("E3").Value = $handoffDate
$wsDe.Range("H3").Value = $epoch
$wsDe.Range("J3").Value = "IsDependency"
$wsDe.Range("K3").Value = "e2e\$newMdName"

$wsDe.Rows.Item(4).Insert()
$wsDe.Range("A4").Value = $png2
$wsDe.Range("B4").Value = ".png"
$wsDe.Range("C4").Value = $ready
$wsDe.Range("D4").Value = $png2target
$wsDe.Range("E4").Value = $handoffDate
$wsDe.Range("H4").Value = $epoch
$wsDe.Range("J4").Value = "IsDependency"
$wsDe.Range("K4").Value = "e2e\$newMdName"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/be33dc629c1990b388306ac50dc5ef87532c4651/e2e/$png1", "", "", $png1)
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/59bd411f2317e6f876b54b36afb791adf1c9a6e8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png1target", "", "", $png1target)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/be33dc629c1990b388306ac50dc5ef87532c4651/e2e/$png2", "", "", $png2)
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/59bd411f2317e6f876b54b36afb791adf1c9a6e8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png2target", "", "", $png2target)

# ---------------------------------------------------------------------------
# Existing "md" hyperlink display text / target needs to reflect the new guid
# ---------------------------------------------------------------------------
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/be33dc629c1990b388306ac50dc5ef87532c4651/e2e/$newMdName", "", "", $newMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/be33dc629c1990b388306ac50dc5ef87532c4651/e2e/$newMdName", "", "", $newMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46cbf783add8b88114af48b7e12de4ed763b940a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", "", "", $zhXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/be33dc629c1990b388306ac50dc5ef87532c4651/e2e/$newMdName", "", "", $newMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/59bd411f2317e6f876b54b36afb791adf1c9a6e8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", "", "", $deXlf)
